$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.051.63'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.853.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '696.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.07'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.852.60'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.20%  '
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("E10").Value = '  +2.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.33'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.463'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("E13").Value = '  +6.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.501.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.859.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '71.124.46'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.75'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.93%  '
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.17'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '496.17'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.10%  '
$ws.Range("E23").Value = '  +1.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.85%  '
$ws.Range("E25").Value = '  +3.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.34'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.72%  '
$ws.Range("E28").Value = '  +2.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.009.83'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.22%  '
$ws.Range("E33").Value = '  +0.35%  '
$ws.Range("E35").Value = '  -0.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.31'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.803.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.98%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  +3.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.38'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +12.85%  '
$ws.Range("E41").Value = '  +1.30%  '
$ws.Range("E42").Value = '  +1.86%  '
$ws.Range("E43").Value = '  +5.99%  '
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '164.61'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.53%  '
$ws.Range("E47").Value = '  +5.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '419.94'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.303'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.85%  '
